# Apply "Ran code for averaged intensities on spiral schemes" edit.
#
# The averaged-intensity table (rows 10-16) is reorganized:
#   - The "Gaussian-Quadrature" row moves up to row 10 (data unchanged).
#   - Three new rows for the spiral sampling schemes are inserted after it
#     (rows 11-13), with freshly computed averaged intensities.
#   - The "NoRotation-tilt60deg" / "Rotation-NoTilt" / "Rotation-60detTilt"
#     rows shift down to rows 14-16 (data unchanged).
#   - The three "HexGrid-..." rows shift down to rows 17-19 (recomputed,
#     matching the previous values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First extend the table: create rows 17-19 by copying the formatting
# (borders/bold/alignment for column A, plain for the rest) from row 16.
$ws.Range("A16:M16").Copy()
$ws.Range("A17:M19").PasteSpecial(-4122)

# --- Row labels (column B) and row index (column A) ---
$labels = @("Gaussian-Quadrature", "Spiral-90deg-10rot-5space", "Spiral-90deg-15rot-5space", "Spiral-90deg-10rot-3space", "NoRotation-tilt60deg", "Rotation-NoTilt", "Rotation-60detTilt", "HexGrid-90degTilt5degRes", "HexGrid-90degTilt22p5degRes", "HexGrid-60degTilt5degRes")

for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = 10 + $i
    $ws.Cells.Item($row, 1).Value = $row - 2
    $ws.Cells.Item($row, 2).Value = $labels[$i]
}

# --- Averaged intensity values (columns C:M) for rows 10-19, in row order ---
$values = @(
    @(0.9997970327128718, 0.9697660149721196, 0.9997385364559017, 0.9997970327128718, 0.9699415614231334, 1.001031999558077, 0.9941468951873086, 0.9697660149721196, 0.9847522757140106, 0.9922746542134413, 0.9890703400515687),
    @(0.9927734059709842, 0.9792377000105068, 0.9947802115617808, 0.9927734059709842, 0.9835745737898508, 1.000151500738855, 0.9941399773397617, 0.9792377000105068, 0.9870089557861438, 0.9898911808785639, 0.9907762282352898),
    @(0.992738379315072, 0.9792811162232559, 0.9947869749823363, 0.992738379315072, 0.9836817641425833, 1.000159340623365, 0.9941451741110581, 0.9792811162232559, 0.9870340456027961, 0.989886212458934, 0.9907987915662785),
    @(0.9927766999458111, 0.9792165617903996, 0.9948272211591419, 0.9927766999458111, 0.9836404378631705, 1.000160864735405, 0.994161525821143, 0.9792165617903996, 0.9870218914747708, 0.989899295710291, 0.9907972185525119),
    @(0.9968679999999996, 0.9312400000000003, 0.9981239999999998, 0.9968679999999996, 0.9341439999999996, 1.053979999999997, 0.9972920000000007, 0.9312400000000003, 0.964682, 0.9807749999999998, 0.9852746666666663),
    @(1, 0.89, 1, 1, 0.89, 1.1, 1, 0.89, 0.9450000000000001, 0.9725, 0.9800000000000001),
    @(0.9968750080000031, 0.9313359026175959, 0.9981250047999984, 0.9968750080000031, 0.9343079706624027, 1.053958451199999, 0.9972916736000018, 0.9313359026175959, 0.9647304537087972, 0.9808027308544002, 0.9853156684800002),
    @(0.9902919605113798, 0.9915977648388012, 0.9911008089428689, 0.9902919605113798, 0.9905667306928275, 0.9900085496758524, 0.9912254889994491, 0.9915977648388012, 0.9913492868908351, 0.9908206237011075, 0.9907985506101964),
    @(0.9893510595467532, 0.9980347295543945, 0.9910885936214049, 0.9893510595467532, 0.9923287446189193, 0.9846493877629974, 0.9902983081668018, 0.9980347295543945, 0.9945616615878997, 0.9919563605673265, 0.9909584705452117),
    @(0.9888161938134855, 1.004513349273281, 0.9886897954349229, 0.9888161938134855, 0.9993333435874322, 0.9801588706975204, 0.9886335434724505, 1.004513349273281, 0.9966015723541017, 0.9927088830837937, 0.9916908493798487)
)

$cols = @("C", "D", "E", "F", "G", "H", "I", "J", "K", "L", "M")

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 10 + $i
    $rowValues = $values[$i]
    for ($c = 0; $c -lt $cols.Length; $c++) {
        $ws.Range($cols[$c] + $row).Value = $rowValues[$c]
    }
}

Write-Host "Averaged intensities for spiral schemes applied"
